# Apply updated simulation results (case with 380 kV done) to Sheet1
# of the active workbook. The diff only changes existing numeric <v>
# values in columns B, D, E, F, G, I, L, M, N for rows 2-25; all other
# cells (A, C, H, J, K, O and the header row) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.119469525981856
$ws.Range("D2").Value = 0.146687881761963
$ws.Range("E2").Value = 0.9711777002425777
$ws.Range("F2").Value = 3.070198749417898
$ws.Range("G2").Value = 0.002445183993796917
$ws.Range("I2").Value = 0.9771647811360893
$ws.Range("L2").Value = 0.6247812563319712
$ws.Range("M2").Value = 0.4085086584114563
$ws.Range("N2").Value = 1.5957035386879
$ws.Range("B3").Value = 1.057582331412334
$ws.Range("D3").Value = 0.1380998653737322
$ws.Range("E3").Value = 0.8453846020535991
$ws.Range("F3").Value = 2.8753568839671
$ws.Range("G3").Value = 0.002456384319351232
$ws.Range("I3").Value = 0.9918490271771176
$ws.Range("L3").Value = 0.5755172734001803
$ws.Range("M3").Value = 0.3815036868041872
$ws.Range("N3").Value = 1.617242365448178
$ws.Range("B4").Value = 1.020129839105323
$ws.Range("D4").Value = 0.1329827212799586
$ws.Range("E4").Value = 0.7682222770407208
$ws.Range("F4").Value = 2.758611451409934
$ws.Range("G4").Value = 0.002463597821144785
$ws.Range("I4").Value = 1.001507900721144
$ws.Range("L4").Value = 0.5456834474105676
$ws.Range("M4").Value = 0.365157621814582
$ws.Range("N4").Value = 1.631148374903223
$ws.Range("B5").Value = 1.005003981176742
$ws.Range("D5").Value = 0.1309354916689358
$ws.Range("E5").Value = 0.7367879808349898
$ws.Range("F5").Value = 2.711742033140354
$ws.Range("G5").Value = 0.002466622437945579
$ws.Range("I5").Value = 1.005604860675973
$ws.Range("L5").Value = 0.5336272144340626
$ws.Range("M5").Value = 0.3585545619886901
$ws.Range("N5").Value = 1.636986596506382
$ws.Range("B6").Value = 1.002500547792209
$ws.Range("D6").Value = 0.1305978152654887
$ws.Range("E6").Value = 0.7315686861136186
$ws.Range("F6").Value = 2.704001395610646
$ws.Range("G6").Value = 0.002467129823544578
$ws.Range("I6").Value = 1.006294854392571
$ws.Range("L6").Value = 0.5316313199290619
$ws.Range("M6").Value = 0.3574616084147237
$ws.Range("N6").Value = 1.637966385973677
$ws.Range("B7").Value = 1.019925295864482
$ws.Range("D7").Value = 0.1329549590873995
$ws.Range("E7").Value = 0.7677983143597942
$ws.Range("F7").Value = 2.757976525713474
$ws.Range("G7").Value = 0.00246363826711925
$ws.Range("I7").Value = 1.001562503250376
$ws.Range("L7").Value = 0.5455204465142458
$ws.Range("M7").Value = 0.365068336746674
$ws.Range("N7").Value = 1.63122641696204
$ws.Range("B8").Value = 1.09801670701313
$ws.Range("D8").Value = 0.1436936599071146
$ws.Range("E8").Value = 0.9277831744423395
$ws.Range("F8").Value = 3.002404980867453
$ws.Range("G8").Value = 0.002448976307862116
$ws.Range("I8").Value = 0.982094084672859
$ws.Range("L8").Value = 0.6077072153707377
$ws.Range("M8").Value = 0.3991479162577747
$ws.Range("N8").Value = 1.602988854606171
$ws.Range("B9").Value = 1.255550155457001
$ws.Range("D9").Value = 0.1660440962434251
$ws.Range("E9").Value = 1.242534771140043
$ws.Range("F9").Value = 3.505637742957617
$ws.Range("G9").Value = 0.002422873114468161
$ws.Range("I9").Value = 0.9490484715132794
$ws.Range("L9").Value = 0.7330854044718649
$ws.Range("M9").Value = 0.4678944528660409
$ws.Range("N9").Value = 1.553012280606843
$ws.Range("B10").Value = 1.374070197888955
$ws.Range("D10").Value = 0.183331478721044
$ws.Range("E10").Value = 1.475029877441244
$ws.Range("F10").Value = 3.891381590038122
$ws.Range("G10").Value = 0.002405280560645509
$ws.Range("I10").Value = 0.9279435212345106
$ws.Range("L10").Value = 0.8275022124985014
$ws.Range("M10").Value = 0.519649319313146
$ws.Range("N10").Value = 1.519574821747156
$ws.Range("B11").Value = 1.428615402928529
$ws.Range("D11").Value = 0.1914019243054952
$ws.Range("E11").Value = 1.581210537968587
$ws.Range("F11").Value = 4.070671029418179
$ws.Range("G11").Value = 0.002397615109945406
$ws.Range("I11").Value = 0.9190415304297588
$ws.Range("L11").Value = 0.8710024941111101
$ws.Range("M11").Value = 0.543482706759832
$ws.Range("N11").Value = 1.505073610259544
$ws.Range("B12").Value = 1.449362578087573
$ws.Range("D12").Value = 0.1944891419676082
$ws.Range("E12").Value = 1.621490083811381
$ws.Range("F12").Value = 4.139138491511005
$ws.Range("G12").Value = 0.002394760423021352
$ws.Range("I12").Value = 0.9157719391111101
$ws.Range("L12").Value = 0.8875578268455513
$ws.Range("M12").Value = 0.5525508906111156
$ws.Range("N12").Value = 1.499684329684882
$ws.Range("B13").Value = 1.444890181919845
$ws.Range("D13").Value = 0.1938228470463343
$ws.Range("E13").Value = 1.61281181206013
$ws.Range("F13").Value = 4.12436682375602
$ws.Range("G13").Value = 0.002395373101594613
$ws.Range("I13").Value = 0.9164715794476237
$ws.Range("L13").Value = 0.8839886006034874
$ws.Range("M13").Value = 0.5505959610951408
$ws.Range("N13").Value = 1.500840471025498
$ws.Range("B14").Value = 1.430320430911195
$ws.Range("D14").Value = 0.1916552795225641
$ws.Range("E14").Value = 1.584522867178208
$ws.Range("F14").Value = 4.07629222295509
$ws.Range("G14").Value = 0.002397379292779954
$ws.Range("I14").Value = 0.918770501645227
$ws.Range("L14").Value = 0.8723628318743692
$ws.Range("M14").Value = 0.5442278822350914
$ws.Range("N14").Value = 1.504628186560979
$ws.Range("B15").Value = 1.421408077007129
$ws.Range("D15").Value = 0.1903316776154327
$ws.Range("E15").Value = 1.567204687790479
$ws.Range("F15").Value = 4.046920743221222
$ws.Range("G15").Value = 0.002398614386484845
$ws.Range("I15").Value = 0.9201918907432045
$ws.Range("L15").Value = 0.8652525971293414
$ws.Range("M15").Value = 0.5403328896003217
$ws.Range("N15").Value = 1.506961556042427
$ws.Range("B16").Value = 1.370518342527475
$ws.Range("D16").Value = 0.1828083304801282
$ws.Range("E16").Value = 1.468100070632715
$ws.Range("F16").Value = 3.879743542804249
$ws.Range("G16").Value = 0.002405788266464819
$ws.Range("I16").Value = 0.9285394260616897
$ws.Range("L16").Value = 0.8246707385091554
$ws.Range("M16").Value = 0.5180976952500771
$ws.Range("N16").Value = 1.520536786923913
$ws.Range("B17").Value = 1.339461512083574
$ws.Range("D17").Value = 0.1782468603051939
$ws.Range("E17").Value = 1.40741736729538
$ws.Range("F17").Value = 3.778180618193858
$ws.Range("G17").Value = 0.002410275310640753
$ws.Range("I17").Value = 0.9338399834103868
$ws.Range("L17").Value = 0.7999184324536373
$ws.Range("M17").Value = 0.5045322805391805
$ws.Range("N17").Value = 1.529046483233898
$ws.Range("B18").Value = 1.321657632719109
$ws.Range("D18").Value = 0.1756425547520735
$ws.Range("E18").Value = 1.372552881709254
$ws.Range("F18").Value = 3.720121555422992
$ws.Range("G18").Value = 0.002412887934928437
$ws.Range("I18").Value = 0.9369544377435162
$ws.Range("L18").Value = 0.7857330635125948
$ws.Range("M18").Value = 0.4967570118426323
$ws.Range("N18").Value = 1.534007832994625
$ws.Range("B19").Value = 1.315639665720198
$ws.Range("D19").Value = 0.1747640585422516
$ws.Range("E19").Value = 1.360754660764428
$ws.Range("F19").Value = 3.700524327394902
$ws.Range("G19").Value = 0.002413777999805343
$ws.Range("I19").Value = 0.9380201972809488
$ws.Range("L19").Value = 0.7809388767155667
$ws.Range("M19").Value = 0.4941290715808293
$ws.Range("N19").Value = 1.535699134546633
$ws.Range("B20").Value = 1.342761428214828
$ws.Range("D20").Value = 0.1787304255388165
$ws.Range("E20").Value = 1.41387306979172
$ws.Range("F20").Value = 3.788955004376078
$ws.Range("G20").Value = 0.002409794369680357
$ws.Range("I20").Value = 0.9332689223091108
$ws.Range("L20").Value = 0.8025480010383603
$ws.Range("M20").Value = 0.5059735177063942
$ws.Range("N20").Value = 1.528133698532557
$ws.Range("B21").Value = 1.434597408151831
$ws.Range("D21").Value = 0.192291090446048
$ws.Range("E21").Value = 1.592829991484336
$ws.Range("F21").Value = 4.090397095340279
$ws.Range("G21").Value = 0.0023967887253483
$ws.Range("I21").Value = 0.918092493071093
$ws.Range("L21").Value = 0.8757753269203761
$ws.Range("M21").Value = 0.5460971636569099
$ws.Range("N21").Value = 1.503512874284567
$ws.Range("B22").Value = 1.49515506366248
$ws.Range("D22").Value = 0.2013358355701484
$ws.Range("E22").Value = 1.710209382725992
$ws.Range("F22").Value = 4.290767995566114
$ws.Range("G22").Value = 0.002388568647999455
$ws.Range("I22").Value = 0.908765296390996
$ws.Range("L22").Value = 0.9241175797859853
$ws.Range("M22").Value = 0.5725715564791898
$ws.Range("N22").Value = 1.488016347669628
$ws.Range("B23").Value = 1.462784595491314
$ws.Range("D23").Value = 0.1964913360558285
$ws.Range("E23").Value = 1.647519491302631
$ws.Range("F23").Value = 4.183509857233787
$ws.Range("G23").Value = 0.002392930410252799
$ws.Range("I23").Value = 0.913688966639377
$ws.Range("L23").Value = 0.8982708742390173
$ws.Range("M23").Value = 0.5584182340767825
$ws.Range("N23").Value = 1.496232737346219
$ws.Range("B24").Value = 1.341269377061451
$ws.Range("D24").Value = 0.178511749289811
$ws.Range("E24").Value = 1.410954376491048
$ws.Range("F24").Value = 3.784082875937798
$ws.Range("G24").Value = 0.002410011700576112
$ws.Range("I24").Value = 0.9335268901215983
$ws.Range("L24").Value = 0.8013590330019156
$ws.Range("M24").Value = 0.5053218607586842
$ws.Range("N24").Value = 1.52854615353651
$ws.Range("B25").Value = 1.212452432079488
$ws.Range("D25").Value = 0.1598513867467091
$ws.Range("E25").Value = 1.157213443027331
$ws.Range("F25").Value = 3.366793769452329
$ws.Range("G25").Value = 0.002429654157485012
$ws.Range("I25").Value = 0.9574341915723004
$ws.Range("L25").Value = 0.6987787217404673
$ws.Range("M25").Value = 0.4490835687916288
$ws.Range("N25").Value = 1.565955533702635

$wb.Save()
